$d = $word.ActiveDocument

$core = "Dates à utiliser pour la Campagne 2018 Persée:  Du 30 octobre au 8 novembre et du 29 novembre au 8 décembre"
$new  = "Dates à utiliser pour la Campagne Cygnus: 10-19 août, 9-18 septembre, 8-17 octobre"

$searchRange = $d.Content
$searchRange.Find.ClearFormatting()

while ($searchRange.Find.Execute($core, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {

    $matchStart = $searchRange.Start
    $matchEnd = $searchRange.End

    # Some occurrences are followed immediately by a trailing space run (different
    # formatting) that belongs to the same sentence and must be absorbed too.
    $peek = $d.Range($matchEnd, $matchEnd + 1).Text
    if ($peek -eq " ") {
        $matchEnd = $matchEnd + 1
    }

    # Replace the whole run sequence [matchStart, matchEnd) with a single plain run
    # (no rPr at all). We purposefully leave the very last original character in
    # place while inserting, then delete it afterwards - fully emptying a
    # paragraph's text before InsertXML causes the host to also replace the
    # paragraph's own <w:pPr> (losing things like <w:sectPr>), which we must not
    # touch here.
    $keepPos = $matchEnd - 1
    $editRange = $d.Range($matchStart, $keepPos)
    $editRange.Text = ""

    $frag = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?><w:wordDocument xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>" + $new + "</w:t></w:r></w:wordDocument>"
    $editRange.InsertXML($frag)

    $leftoverPos = $matchStart + $new.Length
    $leftover = $d.Range($leftoverPos, $leftoverPos + 1)
    $leftover.Text = ""

    $searchRange = $d.Range($leftoverPos, $d.Content.End)
    $searchRange.Find.ClearFormatting()
}
